$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core input change: bump the "Total mass" input (C7). All the
# downstream formulas (C9, C10, C11, C13-C18, E21/E22, B26/C26/E26/F26,
# B27/C27/E27/F27, J16/J17/K16/K17, H15/I15 ...) recalc automatically
# from this single edit. ---
$ws.Range("C7").Value = 1000

# --- New "Keyway" calculation block (rows 32-34) ---
$ws.Range("A32").Value = "Keyway force"
$ws.Range("B32").Formula = "=B26/(20/2)"
$ws.Range("C32").Value = "kN"

$ws.Range("A33").Value = "Keyway area"
$ws.Range("B33").Formula = "=0.02*0.0025*2"

$ws.Range("A34").Value = "Keyway stress"
$ws.Range("B34").Formula = "=(B32)/(B33*1000)"
$ws.Range("C34").Value = "Mpa"

# --- Formatting touch-up: a stray formatted (but empty) cell at B40,
# using a custom 5-decimal numeric format ---
$ws.Range("B40").NumberFormat = "0.00000"

# --- Column B now holds text/labels, widen it to fit ---
$ws.Columns("B").ColumnWidth = 12.71

# --- Selection moved by the author while reviewing the new block ---
$null = $ws.Range("H28").Select()
